# Adding reading and recording
# - Fix header typo on the "Scenario" sheet: workersCount -> workerCount
# - Reformat the "performance" column (C3:C7) on "ProductionCenter" to General
# - Update the active sheet / selection state left behind by the edits

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Scenario")
$ws2 = $wb.Worksheets.Item("ProductionCenter")
$ws3 = $wb.Worksheets.Item("Connection")

# --- Scenario!A2: rename "workersCount" header to "workerCount" ---------
# Pick up the "no wrap" text style already used for id-like cells elsewhere
# in the workbook (ProductionCenter!A3) before writing the corrected text.
[void]$ws2.Range("A3").Copy()
[void]$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("A2").Value = "workerCount"

# --- ProductionCenter!C3:C7: switch number format to General ------------
$ws2.Range("C3:C7").NumberFormat = "General"

# --- Restore per-sheet selections and make Scenario the active tab ------
[void]$ws2.Activate()
[void]$ws2.Range("C3:C7").Select()

[void]$ws3.Activate()
[void]$ws3.Range("E11").Select()

[void]$ws1.Activate()
[void]$ws1.Range("F12").Select()
